$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I8").Value = 'sv'
$ws.Range("J8").Value = 'Statement-opinion'
$ws.Range("I45").Value = 'sd'
$ws.Range("J45").Value = 'Statement-non-opinion'
$ws.Range("I46").Value = 'sv'
$ws.Range("J46").Value = 'Statement-opinion'
$ws.Range("I49").Value = 'aa'
$ws.Range("J49").Value = 'Agree/Accept'
$ws.Range("I70").Value = 'aa'
$ws.Range("J70").Value = 'Agree/Accept'
$ws.Range("I71").Value = '%'
$ws.Range("J71").Value = 'Uninterpretable'
$ws.Range("I76").Value = 'b'
$ws.Range("J76").Value = 'Acknowledge (Backchannel)'
$ws.Range("I83").Value = 'sv'
$ws.Range("J83").Value = 'Statement-opinion'
$ws.Range("I84").Value = 'sd'
$ws.Range("J84").Value = 'Statement-non-opinion'
$ws.Range("I98").Value = 'sd'
$ws.Range("J98").Value = 'Statement-non-opinion'
$ws.Range("I120").Value = 'aa'
$ws.Range("J120").Value = 'Agree/Accept'
$ws.Range("I129").Value = 'sv'
$ws.Range("J129").Value = 'Statement-opinion'
$ws.Range("I137").Value = 'b'
$ws.Range("J137").Value = 'Acknowledge (Backchannel)'
$ws.Range("I142").Value = 'ba'
$ws.Range("J142").Value = 'Appreciation'
$ws.Range("I155").Value = 'aa'
$ws.Range("J155").Value = 'Agree/Accept'
$ws.Range("I164").Value = 'b'
$ws.Range("J164").Value = 'Acknowledge (Backchannel)'
$ws.Range("I181").Value = 'sd'
$ws.Range("J181").Value = 'Statement-non-opinion'
$ws.Range("I182").Value = 'sd'
$ws.Range("J182").Value = 'Statement-non-opinion'
$ws.Range("I185").Value = 'b'
$ws.Range("J185").Value = 'Acknowledge (Backchannel)'
$ws.Range("I187").Value = 'sv'
$ws.Range("J187").Value = 'Statement-opinion'
$ws.Range("I189").Value = 'b'
$ws.Range("J189").Value = 'Acknowledge (Backchannel)'
$ws.Range("I221").Value = 'b'
$ws.Range("J221").Value = 'Acknowledge (Backchannel)'
$ws.Range("I226").Value = 'sd'
$ws.Range("J226").Value = 'Statement-non-opinion'
$ws.Range("I228").Value = 'aa'
$ws.Range("J228").Value = 'Agree/Accept'
$ws.Range("I232").Value = 'sv'
$ws.Range("J232").Value = 'Statement-opinion'
$ws.Range("I242").Value = 'b'
$ws.Range("J242").Value = 'Acknowledge (Backchannel)'
$ws.Range("I245").Value = 'b'
$ws.Range("J245").Value = 'Acknowledge (Backchannel)'
$ws.Range("I246").Value = 'sv'
$ws.Range("J246").Value = 'Statement-opinion'
$ws.Range("I279").Value = 'sd'
$ws.Range("J279").Value = 'Statement-non-opinion'
$ws.Range("I280").Value = 'sd'
$ws.Range("J280").Value = 'Statement-non-opinion'
$ws.Range("I287").Value = 'b'
$ws.Range("J287").Value = 'Acknowledge (Backchannel)'
$ws.Range("I291").Value = 'sd'
$ws.Range("J291").Value = 'Statement-non-opinion'
$ws.Range("I293").Value = 'b'
$ws.Range("J293").Value = 'Acknowledge (Backchannel)'
$ws.Range("I297").Value = '%'
$ws.Range("J297").Value = 'Uninterpretable'
$ws.Range("I299").Value = 'aa'
$ws.Range("J299").Value = 'Agree/Accept'
$ws.Range("I300").Value = 'aa'
$ws.Range("J300").Value = 'Agree/Accept'
$ws.Range("I302").Value = 'sd'
$ws.Range("J302").Value = 'Statement-non-opinion'
$ws.Range("I310").Value = 'sd'
$ws.Range("J310").Value = 'Statement-non-opinion'
$ws.Range("I321").Value = 'sv'
$ws.Range("J321").Value = 'Statement-opinion'
$ws.Range("I322").Value = 'aa'
$ws.Range("J322").Value = 'Agree/Accept'
$ws.Range("I345").Value = 'sd'
$ws.Range("J345").Value = 'Statement-non-opinion'
$ws.Range("I362").Value = 'aa'
$ws.Range("J362").Value = 'Agree/Accept'
$ws.Range("I373").Value = 'qy'
$ws.Range("J373").Value = 'Yes-No-Question'
$ws.Range("I390").Value = 'ba'
$ws.Range("J390").Value = 'Appreciation'
$ws.Range("I392").Value = 'sv'
$ws.Range("J392").Value = 'Statement-opinion'
$ws.Range("I396").Value = '%'
$ws.Range("J396").Value = 'Uninterpretable'
$ws.Range("I406").Value = '%'
$ws.Range("J406").Value = 'Uninterpretable'
$ws.Range("I416").Value = 'sd'
$ws.Range("J416").Value = 'Statement-non-opinion'
$ws.Range("I422").Value = 'sv'
$ws.Range("J422").Value = 'Statement-opinion'
$ws.Range("I428").Value = 'sd'
$ws.Range("J428").Value = 'Statement-non-opinion'
$ws.Range("I433").Value = 'sv'
$ws.Range("J433").Value = 'Statement-opinion'
$ws.Range("I456").Value = 'sd'
$ws.Range("J456").Value = 'Statement-non-opinion'
$ws.Range("I469").Value = 'sd'
$ws.Range("J469").Value = 'Statement-non-opinion'
$ws.Range("I485").Value = 'sd'
$ws.Range("J485").Value = 'Statement-non-opinion'
$ws.Range("I499").Value = 'sd'
$ws.Range("J499").Value = 'Statement-non-opinion'
$ws.Range("I503").Value = 'sd'
$ws.Range("J503").Value = 'Statement-non-opinion'
$ws.Range("I508").Value = 'sd'
$ws.Range("J508").Value = 'Statement-non-opinion'
$ws.Range("I513").Value = 'ba'
$ws.Range("J513").Value = 'Appreciation'
$ws.Range("I514").Value = 'sd'
$ws.Range("J514").Value = 'Statement-non-opinion'
$ws.Range("I523").Value = 'b'
$ws.Range("J523").Value = 'Acknowledge (Backchannel)'
$ws.Range("I525").Value = 'aa'
$ws.Range("J525").Value = 'Agree/Accept'
$ws.Range("I527").Value = 'sd'
$ws.Range("J527").Value = 'Statement-non-opinion'
$ws.Range("I529").Value = 'sd'
$ws.Range("J529").Value = 'Statement-non-opinion'
$ws.Range("I536").Value = 'sv'
$ws.Range("J536").Value = 'Statement-opinion'
$ws.Range("I542").Value = 'sd'
$ws.Range("J542").Value = 'Statement-non-opinion'
$ws.Range("I548").Value = 'sd'
$ws.Range("J548").Value = 'Statement-non-opinion'
$ws.Range("I561").Value = 'sd'
$ws.Range("J561").Value = 'Statement-non-opinion'
